$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Range("H15").Value = 2627.4546
$ws.Range("I15").Value = 2627.4546
$ws.Range("K15").Value = 7882.3638
$ws.Range("M15").Value = -7713.3638
# Row 19
$ws.Range("H19").Value = 1331.2222
$ws.Range("I19").Value = 1214.1428
$ws.Range("K19").Value = 1214.1428
$ws.Range("M19").Value = -1039.1428
$ws.Range("N19").ClearContents()
# Row 37
$ws.Range("H37").Value = 550
$ws.Range("I37").Value = 100
$ws.Range("J37").Value = 1000
$ws.Range("K37").Value = 300
$ws.Range("L37").Value = 3000
$ws.Range("M37").Value = -174
$ws.Range("N37").Value = -3252
# Row 40
$ws.Range("H40").Value = 1240.4117
$ws.Range("I40").Value = 1249.0714
$ws.Range("J40").Value = 1200
$ws.Range("K40").Value = 1249.0714
$ws.Range("L40").Value = 1200
$ws.Range("M40").Value = -1074.0714
$ws.Range("N40").Value = -1550
# Row 55
$ws.Range("H55").Value = 756.8
$ws.Range("J55").Value = 571
$ws.Range("L55").Value = 571
$ws.Range("N55").Value = -999
# Row 70
$ws.Range("H70").Value = 4667.6206
$ws.Range("I70").Value = 3947.2144
$ws.Range("J70").Value = 5340
$ws.Range("K70").Value = 11841.6432
$ws.Range("L70").Value = 16020
$ws.Range("M70").Value = -11571.6432
$ws.Range("N70").Value = -16560
# Row 73
$ws.Range("H73").Value = 4667.6206
$ws.Range("I73").Value = 3947.2144
$ws.Range("J73").Value = 5340
$ws.Range("K73").Value = 11841.6432
$ws.Range("L73").Value = 16020
$ws.Range("M73").Value = -10905.6432
$ws.Range("N73").Value = -17892
# Row 74
$ws.Range("H74").Value = 6166
$ws.Range("I74").Value = 5199.2
$ws.Range("K74").Value = 5199.2
$ws.Range("M74").Value = -4263.2
$ws.Range("N74").ClearContents()
# Row 77
$ws.Range("H77").Value = 6166
$ws.Range("I77").Value = 5199.2
$ws.Range("K77").Value = 25996
$ws.Range("M77").Value = -21316
$ws.Range("N77").ClearContents()
# Row 80
$ws.Range("H80").Value = 551.375
$ws.Range("I80").Value = 412.7
$ws.Range("J80").Value = 650.4286
$ws.Range("K80").Value = 1238.1
$ws.Range("L80").Value = 1951.2858
$ws.Range("M80").Value = -240.0999999999999
$ws.Range("N80").Value = -3947.2858
# Row 83
$ws.Range("H83").Value = 551.375
$ws.Range("I83").Value = 412.7
$ws.Range("J83").Value = 650.4286
$ws.Range("K83").Value = 3714.3
$ws.Range("L83").Value = 5853.8574
$ws.Range("M83").Value = 1277.7
$ws.Range("N83").Value = -15837.8574
# Row 118
$ws.Range("H118").Value = 269
$ws.Range("J118").Value = 149
$ws.Range("L118").Value = 447
$ws.Range("N118").Value = -3761
# Row 132
$ws.Range("H132").Value = 3476.2593
$ws.Range("I132").Value = 3599.0833
$ws.Range("J132").Value = 2493.6667
$ws.Range("K132").Value = 10797.2499
$ws.Range("L132").Value = 7481.000100000001
$ws.Range("M132").Value = -8267.249899999999
$ws.Range("N132").Value = -12541.0001

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 43
$ws.Range("H43").Value = 29471.75
$ws.Range("J43").Value = 29471.75
$ws.Range("L43").Value = 29471.75
$ws.Range("N43").Value = -30097.75
# Row 62
$ws.Range("H62").Value = 24300
$ws.Range("I62").Value = 24300
$ws.Range("K62").Value = 24300
$ws.Range("M62").Value = -23676
# Row 65
$ws.Range("H65").Value = 24300
$ws.Range("I65").Value = 24300
$ws.Range("K65").Value = 72900
$ws.Range("M65").Value = -69780
# Row 101
$ws.Range("H101").Value = 24499.75
$ws.Range("J101").Value = 24499.75
$ws.Range("L101").Value = 24499.75
$ws.Range("N101").Value = -30989.75
# Row 122
$ws.Range("H122").Value = 3068.7693
$ws.Range("I122").Value = 3207.8333
$ws.Range("K122").Value = 9623.499899999999
$ws.Range("M122").Value = -7173.499899999999
$ws.Range("N122").ClearContents()
# Row 132
$ws.Range("H132").Value = 2162.6978
$ws.Range("I132").Value = 2098.65
$ws.Range("K132").Value = 6295.950000000001
$ws.Range("M132").Value = -3765.950000000001
$ws.Range("N132").ClearContents()

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 59
$ws.Range("H59").Value = 60000
$ws.Range("I59").Value = 60000
$ws.Range("K59").Value = 60000
$ws.Range("M59").Value = -59153

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2783.3684
$ws.Range("I31").Value = 1422.8572
$ws.Range("K31").Value = 1422.8572
$ws.Range("M31").Value = -1127.8572
$ws.Range("N31").ClearContents()
# Row 34
$ws.Range("H34").Value = 2783.3684
$ws.Range("I34").Value = 1422.8572
$ws.Range("K34").Value = 1422.8572
$ws.Range("M34").Value = -1220.8572
$ws.Range("N34").ClearContents()
# Row 37
$ws.Range("H37").Value = 24994.5
$ws.Range("J37").Value = 24992.666
$ws.Range("L37").Value = 24992.666
$ws.Range("N37").Value = -25206.666
# Row 41
$ws.Range("H41").Value = 25077.766
$ws.Range("I41").Value = 9472
$ws.Range("J41").Value = 38949.555
$ws.Range("K41").Value = 9472
$ws.Range("L41").Value = 38949.555
$ws.Range("M41").Value = -9044
$ws.Range("N41").Value = -39805.555
# Row 58
$ws.Range("H58").Value = 2474
$ws.Range("I58").Value = 1538.6
$ws.Range("K58").Value = 1538.6
$ws.Range("M58").Value = -1335.6
$ws.Range("N58").ClearContents()
# Row 105
$ws.Range("H105").Value = 3380.25
$ws.Range("I105").Value = 3170
$ws.Range("K105").Value = 3170
$ws.Range("M105").Value = -1423
$ws.Range("N105").ClearContents()
# Row 122
$ws.Range("H122").Value = 3409.5454
$ws.Range("I122").Value = 3925.8667
$ws.Range("K122").Value = 11777.6001
$ws.Range("M122").Value = -9327.6001
$ws.Range("N122").ClearContents()
# Row 132
$ws.Range("H132").Value = 1192
$ws.Range("I132").Value = 1192
$ws.Range("K132").Value = 3576
$ws.Range("M132").Value = -1046
# Row 134
$ws.Range("H134").Value = 1903.4073
$ws.Range("I134").Value = 1973.6086
$ws.Range("K134").Value = 5920.825800000001
$ws.Range("M134").Value = -3385.825800000001
$ws.Range("N134").ClearContents()
# Row 136
$ws.Range("H136").Value = 2474
$ws.Range("I136").Value = 1538.6
$ws.Range("K136").Value = 4615.799999999999
$ws.Range("M136").Value = -2065.799999999999
$ws.Range("N136").ClearContents()

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 38
$ws.Range("H38").Value = 222
$ws.Range("J38").Value = 86
$ws.Range("L38").Value = 258
$ws.Range("N38").Value = -952
# Row 39
$ws.Range("H39").Value = 5417.0835
$ws.Range("J39").Value = 5636.5454
$ws.Range("L39").Value = 16909.6362
$ws.Range("N39").Value = -17497.6362
# Row 60
$ws.Range("H60").Value = 1750
$ws.Range("I60").Value = 1000
$ws.Range("K60").Value = 3000
$ws.Range("M60").Value = -2749
$ws.Range("N60").ClearContents()
# Row 98
$ws.Range("H98").Value = 1221
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 1221
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 3663
$ws.Range("N98").Value = -6659
$ws.Range("M98").ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 57
$ws.Range("H57").Value = 11193
$ws.Range("I57").Value = 5000
$ws.Range("J57").Value = 14289.5
$ws.Range("K57").Value = 5000
$ws.Range("L57").Value = 14289.5
$ws.Range("M57").Value = -4180
$ws.Range("N57").Value = -15929.5
# Row 102
$ws.Range("H102").Value = 1931
$ws.Range("I102").Value = 1807.4615
$ws.Range("K102").Value = 1807.4615
$ws.Range("M102").Value = -185.4614999999999
$ws.Range("N102").ClearContents()
# Row 113
$ws.Range("H113").Value = 1878.4286
$ws.Range("I113").Value = 1037.5
$ws.Range("K113").Value = 1037.5
$ws.Range("M113").Value = 1132.5
$ws.Range("N113").ClearContents()
# Row 132
$ws.Range("H132").Value = 3003.7307
$ws.Range("I132").Value = 2716.652
$ws.Range("K132").Value = 8149.956
$ws.Range("M132").Value = -5619.956
$ws.Range("N132").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 68
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("M68").ClearContents()
$ws.Range("N68").ClearContents()
# Row 71
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("M71").ClearContents()
$ws.Range("N71").ClearContents()
# Row 94
$ws.Range("H94").Value = 46799.4
$ws.Range("J94").Value = 46799.4
$ws.Range("L94").Value = 46799.4
$ws.Range("N94").Value = -48151.4
# Row 100
$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("M100").ClearContents()
$ws.Range("N100").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 53
$ws.Range("H53").Value = 30000
$ws.Range("I53").Value = 30000
$ws.Range("K53").Value = 30000
$ws.Range("M53").Value = -29393
# Row 136
$ws.Range("H136").Value = 8133.095
$ws.Range("I136").Value = 4649.6665
$ws.Range("K136").Value = 13948.9995
$ws.Range("M136").Value = -11398.9995
$ws.Range("N136").ClearContents()

Write-Host "Applied all Phantom_Profits updates"